$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.536.21"
$ws.Range("E2").Value = "  -5.88%  "
$ws.Range("D3").Value = "3.046.81"
$ws.Range("E3").Value = "  -6.42%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "538.99"
$ws.Range("E5").Value = "  -7.40%  "
$ws.Range("D6").Value = "133.32"
$ws.Range("E6").Value = "  -13.55%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.053.18"
$ws.Range("E8").Value = "  -6.02%  "
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  -7.16%  "
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -13.46%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -6.42%  "
$ws.Range("D13").Value = "34.54"
$ws.Range("E13").Value = "  -9.24%  "
$ws.Range("D14").Value = "0.0000212"
$ws.Range("E14").Value = "  -10.05%  "
$ws.Range("D15").Value = "3.565.13"
$ws.Range("E15").Value = "  -5.73%  "
$ws.Range("D16").Value = "62.691.91"
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").Value = "3.074.30"
$ws.Range("E18").Value = "  -5.62%  "
$ws.Range("D19").Value = "6.60"
$ws.Range("E19").Value = "  -7.78%  "
$ws.Range("D20").Value = "477.15"
$ws.Range("E20").Value = "  -13.99%  "
$ws.Range("D21").Value = "13.32"
$ws.Range("E21").Value = "  -8.30%  "
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -5.85%  "
$ws.Range("D23").Value = "7.11"
$ws.Range("E23").Value = "  -9.65%  "
$ws.Range("D24").Value = "77.79"
$ws.Range("E24").Value = "  -4.90%  "
$ws.Range("D25").Value = "12.06"
$ws.Range("E25").Value = "  -12.08%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -9.68%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "8.25"
$ws.Range("E28").Value = "  -11.21%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "26.08"
$ws.Range("E30").Value = "  -6.26%  "
$ws.Range("D31").Value = "1.91"
$ws.Range("E31").Value = "  -15.62%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "60.40"
$ws.Range("E32").Value = "  +9.33%  "
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").Value = "1.10"
$ws.Range("E33").Value = "  -6.69%  "
$ws.Range("D34").Value = "2.43"
$ws.Range("E34").Value = "  -12.24%  "
$ws.Range("D35").Value = "492.34"
$ws.Range("E35").Value = "  -13.18%  "
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -8.03%  "
$ws.Range("D37").Value = "4.99"
$ws.Range("E37").Value = "  -12.84%  "
$ws.Range("D38").Value = "3.119.00"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").Value = "0.0389"
$ws.Range("E39").Value = "  -14.46%  "
$ws.Range("D40").Value = "0.0785"
$ws.Range("E40").Value = "  -9.24%  "
$ws.Range("D41").Value = "0.116"
$ws.Range("E41").Value = "  -10.21%  "
$ws.Range("D42").Value = "7.98"
$ws.Range("E42").Value = "  -7.90%  "
$ws.Range("D43").Value = "2.53"
$ws.Range("E43").Value = "  -16.60%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").Value = "  -10.59%  "
$ws.Range("D46").Value = "24.65"
$ws.Range("E46").Value = "  -7.02%  "
$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  -13.58%  "
$ws.Range("D48").Value = "117.52"
$ws.Range("E48").Value = "  -6.92%  "
$ws.Range("D49").Value = "0.107"
$ws.Range("E49").Value = "  -5.65%  "
$ws.Range("D50").Value = "0.0₃0497"
$ws.Range("E50").Value = "  -11.57%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  +18.24%  "
